$wb = $excel.ActiveWorkbook

# --- Sheet: Function Index ---
$ws1 = $wb.Worksheets.Item("Function Index")
$ws1.Range("D2").Value = "setAspect, nowMs, init, evaluate"
$ws1.Range("D3").Value = "nowMs, isOccupied, evaluateControllerLogic, update, setAspect, isHealthy"
$ws1.Range("E4").Value = $false
$ws1.Range("E5").Value = $false
$ws1.Range("D7").Value = "readRawClear, configure"
$ws1.Range("E7").Value = $false
$ws1.Range("D9").Value = "else, pinMode"
$ws1.Range("D10").Value = "else, digitalWrite"

# --- Sheet: Call Graph ---
$ws2 = $wb.Worksheets.Item("Call Graph")
$ws2.Range("B2").Value = "setAspect, nowMs, init, evaluate"
$ws2.Range("B3").Value = "nowMs, isOccupied, evaluateControllerLogic, update, setAspect, isHealthy"
$ws2.Range("B7").Value = "readRawClear, configure"
$ws2.Range("B9").Value = "else, pinMode"
$ws2.Range("B10").Value = "else, digitalWrite"

# --- Sheet: File Summaries ---
$ws3 = $wb.Worksheets.Item("File Summaries")
$ws3.Range("D3").Value = $true
$ws3.Range("D4").Value = $true
$ws3.Range("F10").Value = "src\logic\ControllerHelpers.cpp, src\logic\Interlocking.cpp"

# --- Sheet: Class Roles ---
# Reorder rows 3-6 (A3:B6) to: MockGpio/HARDWARE, SignalHead/MIXED, BlockController/MIXED, ArduinoGpio/HARDWARE
$ws4 = $wb.Worksheets.Item("Class Roles")
$ws4.Range("A3").Value = "MockGpio"
$ws4.Range("B3").Value = "HARDWARE"
$ws4.Range("A4").Value = "SignalHead"
$ws4.Range("B4").Value = "MIXED"
$ws4.Range("A5").Value = "BlockController"
$ws4.Range("B5").Value = "MIXED"
$ws4.Range("A6").Value = "ArduinoGpio"
$ws4.Range("B6").Value = "HARDWARE"
